$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-38 down to 24-39
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly data record
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44467
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112009
$ws.Range("G23").Value = "Acelga"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 1200
$ws.Range("M23").Value = 1100
$ws.Range("N23").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 367
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = "Hortaliza"
